$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) ----
$ws.Range("A1").Value = "Tarih"
$ws.Range("B1").Value = "Kullanıcı"
$ws.Range("C1").Value = "Ruh Hali"
$ws.Range("D1").Value = "Müzik"
$ws.Range("E1").Value = "Aktivite"
$ws.Range("F1").Value = "Mesaj"

# ---- Row 2 ----
$ws.Range("A2").Value = "2025-11-22 19:59:10"
$ws.Range("B2").Value = "merve"
$ws.Range("C2").Value = "Mutlu"
$ws.Range("D2").Value = "Happy - Pharrell"
$ws.Range("E2").Value = "Dans etmek"
$ws.Range("F2").Value = ""

# ---- Row 3 ----
$ws.Range("A3").Value = "2025-11-22 19:59:19"
$ws.Range("B3").Value = "ibrahim"
$ws.Range("C3").Value = "Mutlu"
$ws.Range("D3").Value = "Happy - Pharrell"
$ws.Range("E3").Value = "Dans etmek"
$ws.Range("F3").Value = ""

# ---- Row 4 ----
$ws.Range("A4").Value = "2025-11-22 19:59:28"
$ws.Range("B4").Value = "hanife"
$ws.Range("C4").Value = "Mutlu"
$ws.Range("D4").Value = "Happy - Pharrell"
$ws.Range("E4").Value = "Dans etmek"
$ws.Range("F4").Value = ""

# ---- Row 5 ----
$ws.Range("A5").Value = "2025-11-22 20:02:00"
$ws.Range("B5").Value = ""
$ws.Range("C5").Value = "Yorgun"
$ws.Range("D5").Value = "Easy On Me - Adele"
$ws.Range("E5").Value = "Kısa uyku"
$ws.Range("F5").Value = ""

# ---- Row 6 (real date/time value, written twice with two different
#      NumberFormat casings to reproduce the two numFmt entries) ----
$ws.Range("A6").Value = 45983.83710914352
$ws.Range("A6").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("A6").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B6").Value = "ahmet"
$ws.Range("C6").Value = "Motivasyonlu"
$ws.Range("D6").Value = "Stronger - Kanye"
$ws.Range("E6").Value = "Hedef belirle"
$ws.Range("F6").Value = "Devam et, harika işler başarabilirsin 💪"

# ---- Row 7 (real date/time value) ----
$ws.Range("A7").Value = 45983.83731133936
$ws.Range("A7").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B7").Value = "ahmet"
$ws.Range("C7").Value = "Motivasyonlu"
$ws.Range("D7").Value = "Eye of the Tiger - Survivor"
$ws.Range("E7").Value = "Hedef belirle"
$ws.Range("F7").Value = "Devam et, harika işler başarabilirsin 💪"
